$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values would otherwise be auto-converted to numbers
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply cell value updates
$ws.Range("D2").Value = '29.194.73'
$ws.Range("E2").Value = '  +0.18%  '
$ws.Range("D3").Value = '1.836.46'
$ws.Range("E3").Value = '  -0.06%  '
$ws.Range("D4").Value = '0.9997'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '241.22'
$ws.Range("E5").Value = '  +0.38%  '
$ws.Range("D6").Value = '0.6642'
$ws.Range("E6").Value = '  -2.96%  '
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '0.07363'
$ws.Range("E8").Value = '  -0.65%  '
$ws.Range("D9").Value = '0.2925'
$ws.Range("E9").Value = '  -2.22%  '
$ws.Range("D10").Value = '22.65'
$ws.Range("E10").Value = '  -2.01%  '
$ws.Range("E11").Value = '  +1.01%  '
$ws.Range("D12").Value = '1.833.19'
$ws.Range("E12").Value = '  -0.17%  '
$ws.Range("D13").Value = '4.981'
$ws.Range("E13").Value = '  -1.47%  '
$ws.Range("D14").Value = '0.6684'
$ws.Range("E14").Value = '  -1.68%  '
$ws.Range("D15").Value = '82.68'
$ws.Range("E15").Value = '  -5.45%  '
$ws.Range("D16").Value = '6.097'
$ws.Range("E16").Value = '  -0.89%  '
$ws.Range("D17").Value = '29.182.08'
$ws.Range("E17").Value = '  +0.12%  '
$ws.Range("D18").Value = '0.000008280'
$ws.Range("D19").Value = '225.50'
$ws.Range("E19").Value = '  -1.45%  '
$ws.Range("E20").Value = '  -0.67%  '
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  +0.14%  '
$ws.Range("D22").Value = '7.122'
$ws.Range("E22").Value = '  -2.94%  '
$ws.Range("E23").Value = '  +0.07%  '
$ws.Range("D24").Value = '160.49'
$ws.Range("E24").Value = '  +0.50%  '
$ws.Range("D25").Value = '8.620'
$ws.Range("E25").Value = '  -0.94%  '
$ws.Range("D26").Value = '0.1392'
$ws.Range("E26").Value = '  -3.24%  '
$ws.Range("E27").Value = '  -0.44%  '
$ws.Range("D28").Value = '1.511'
$ws.Range("E28").Value = '  -0.05%  '
$ws.Range("D29").Value = '4.107'
$ws.Range("E29").Value = '  -3.58%  '
$ws.Range("D30").Value = '4.040'
$ws.Range("E30").Value = '  -2.34%  '
$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").Value = '0.05343'
$ws.Range("E31").Value = '  +1.55%  '
$ws.Range("B32").Value = 'Toncoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D32").Value = '1.178'
$ws.Range("E32").Value = '  -1.26%  '
$ws.Range("D33").Value = '1.871'
$ws.Range("E33").Value = '  +1.20%  '
$ws.Range("D34").Value = '0.7535'
$ws.Range("E34").Value = '  -0.14%  '
$ws.Range("D35").Value = '1.130'
$ws.Range("E35").Value = '  -0.16%  '
$ws.Range("D36").Value = '2.674'
$ws.Range("E36").Value = '  -0.30%  '
$ws.Range("D37").Value = '1.289.11'
$ws.Range("E37").Value = '  -0.17%  '
$ws.Range("D38").Value = '0.01794'
$ws.Range("E38").Value = '  -1.63%  '
$ws.Range("D39").Value = '2.721'
$ws.Range("E39").Value = '  +0.04%  '
$ws.Range("D40").Value = '0.9223'
$ws.Range("E40").Value = '  -1.65%  '
$ws.Range("D41").Value = '0.09012'
$ws.Range("E41").Value = '  +20.72%  '
$ws.Range("D42").Value = '5.954'
$ws.Range("E42").Value = '  +0.40%  '
$ws.Range("D43").Value = '1.005'
$ws.Range("E43").Value = '  +0.53%  '
$ws.Range("D44").Value = '102.24'
$ws.Range("E44").Value = '  -2.20%  '
$ws.Range("D45").Value = '1.977.61'
$ws.Range("E45").Value = '  -0.19%  '
$ws.Range("D46").Value = '0.5164'
$ws.Range("E46").Value = '  -0.55%  '
$ws.Range("D47").Value = '1.767'
$ws.Range("E47").Value = '  +0.28%  '
$ws.Range("E48").Value = '  -2.51%  '
$ws.Range("D49").Value = '63.00'
$ws.Range("E49").Value = '  -2.41%  '
$ws.Range("D50").Value = '0.05930'
$ws.Range("D51").Value = '9.022'
$ws.Range("E51").Value = '  -4.47%  '
